{"js": "// 1. Replace the author name \"Yannis Plaschko\" with \"Maximilian Meier\"\n//    on the title-page \"Author: \" line (leave the later \"Supplier\" list,\n//    which still legitimately lists Yannis Plaschko, untouched).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet authorParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Author:\") !== -1) {\n    authorParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (authorParagraph) {\n  const nameHits = authorParagraph.search(\"Yannis Plaschko\", { matchCase: true });\n  nameHits.load(\"items\");\n  await context.sync();\n  if (nameHits.items.length > 0) {\n    nameHits.items[0].insertText(\"Maximilian Meier\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 2. Append a new row to the \"Version / Date / Author / Comment\" history\n//    table (the first table in the document) recording the format fix.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst historyTable = tables.items[0];\nhistoryTable.addRows(Word.InsertLocation.end, 1, [\n  [\"1.1\", \"05.05.2022\", \"Maximilian Meier\", \"Document format fix\"]\n]);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Replace the author name \"Yannis Plaschko\" with \"Maximilian Meier\"\n#    on the title-page \"Author: \" line (leave the later \"Supplier\" list,\n#    which still legitimately lists Yannis Plaschko, untouched).\n$authorPara = $null\nforeach ($p in $d.Paragraphs) {\n  if ($p.Range.Text -like \"Author:*\") {\n    $authorPara = $p\n    break\n  }\n}\n\nif ($authorPara -ne $null) {\n  $rng = $authorPara.Range\n  $rng.Find.Execute(\"Yannis Plaschko\", $false, $false, $false, $false, $false, $true, 1, $false, \"Maximilian Meier\", 2)\n}\n\n# 2. Append a new row to the \"Version / Date / Author / Comment\" history\n#    table (the first table in the document) recording the format fix.\n$t = $d.Tables.Item(1)\n$newRow = $t.Rows.Add()\n$newRow.Cells.Item(1).Range.Text = \"1.1\"\n$newRow.Cells.Item(2).Range.Text = \"05.05.2022\"\n$newRow.Cells.Item(3).Range.Text = \"Maximilian Meier\"\n$newRow.Cells.Item(4).Range.Text = \"Document format fix\"\n"}
